# Sound Asset List.xlsx - add more asset-list rows, a closing note, and
# related formatting (column width, row heights, wrap text, header border,
# page orientation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data -------------------------------------------------------------
# Columns: A=Event, B=Description, C=Categorization, D=Assets Required, E=Status
$rows = @(
    @{ Row = 2;  A = "Initial game startup";
                 B = "Main menu theme (think Mass Effect but less music/less serious)";
                 C = "Ambience";
                 D = "Hollow, drawn out noises like you would hear in space movies";
                 E = "Recording"; Height = 47.25 },
    @{ Row = 3;  A = "Player presses enter to begin game";
                 B = "Echoing ding-like sound; indicates that the player pressed a key to start";
                 C = "Interface";
                 D = "A solid dinging noise";
                 E = "Recording"; Height = 47.25 },
    @{ Row = 4;  A = "Player enters the main game/leaves initial main screen";
                 B = "Entry to game, robot explaining what happened to the astronauts";
                 C = "Dialog";
                 D = "Voiceover (robotized)";
                 E = "Recording"; Height = 27.75; NoWrapD = $true },
    @{ Row = 5;  A = "Player is in-game (constant background ambience)";
                 B = "Space traversal music, general game background noises/ambience";
                 C = "Ambience";
                 D = "More hollow, echoing space sounds";
                 E = "Recording"; Height = 29.25 },
    @{ Row = 6;  A = "Player changes the tilt position of the ship (Q/E)";
                 B = "Ship rotates side to side, release a ping adjustment sound";
                 C = "Sound Effect";
                 D = "A pinging sound. Like a submarine but with less echo";
                 E = "Recording"; Height = 45 },
    @{ Row = 7;  A = "Player moves ship forwards";
                 B = "Ship moves forwards, deep whooshing noise that gets higher in tone as time goes on";
                 C = "Sound Effect";
                 D = "Whoosh effect, higher tones as it continues";
                 E = "Recording"; Height = 30 },
    @{ Row = 8;  A = "Player moves ship backwards";
                 B = "Ship moves backwards, deep whoosing noise that gets lower in tone as time goes on";
                 C = "Sound Effect";
                 D = "Whoosh effect, deeper tones as it continues";
                 E = "Recording"; Height = 30 },
    @{ Row = 9;  A = "Player turns ship left or right";
                 B = "Ship turns to either side, lower whooshing tone so as not to interfere with forward or backward sounds";
                 C = "Sound Effect";
                 D = "Another whoosh, lower register";
                 E = "Recording"; Height = 30 },
    @{ Row = 10; A = "Player hits escape to enter the pause menu";
                 B = "High-tech sound, like a TV screen turning on; electric sound, `"futuristic user interface`" sound";
                 C = "Interface";
                 D = "high-pitched whir, electric buzz";
                 E = "Recording"; Height = 30 },
    @{ Row = 11; A = "Player approaches one of the astronauts";
                 B = "Astronaut sends an SOS signal over their radio, gets louder as ship approaches";
                 C = "Sound Effect/Dialog";
                 D = "distress signal, voiceover";
                 E = "Recording"; Height = $null }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    if ($r.Height) {
        $ws.Rows.Item($n).RowHeight = $r.Height
    }
}

# Row 12: only the event/column-A note about picking up an astronaut.
$ws.Range("A12").Value = "Player picks up an astronaut (runs into them with the ship)"

# Closing note, further down the sheet.
$ws.Range("B19").Value = "*Sound List subject to additions or changes*"

# --- Formatting -------------------------------------------------------------

# Wrap text on the "Assets Required" column for every data row except row 4
# (which keeps its original, unwrapped style).
$wrapRange = $ws.Range("D2,D3,D5,D6,D7,D8,D9,D10,D11")
$wrapRange.WrapText = $true

# Center-align the closing note.
$ws.Range("B19").HorizontalAlignment = -4108

# Widen the Description column and select the new note cell, matching the
# final view state of the workbook.
$ws.Columns.Item(2).ColumnWidth = 91.6
$ws.Range("B12").Select()

# Add a thin bottom border under the header row (Event/Description/...).
$headerBorder = $ws.Range("A1:E1").Borders.Item(9)
$headerBorder.LineStyle = 1
$headerBorder.Weight = 2

# Print in portrait orientation.
$ws.PageSetup.Orientation = 1
